$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "_old" suffixed headers (columns A-J) to "_FV2404"
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value()
    $cell.Value = ($val -replace "_old$", "_FV2404")
}

# Rename the "_new" suffixed headers (columns L-U) to "_FV2410"
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value()
    $cell.Value = ($val -replace "_new$", "_FV2410")
}

# Turn the data range into an Excel Table ("Table1") covering A1:U93
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U93"), 0, 1)
$tbl.Name = "Table1"

# Freeze the header row (split below row 1)
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
